$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 13467
$ws.Range("I6").Value = 14056.846
$ws.Range("K6").Value = 42170.538
$ws.Range("M6").Value = -42058.538

$ws.Range("H15").Value = 731.63794
$ws.Range("I15").Value = 731.63794
$ws.Range("K15").Value = 2194.91382
$ws.Range("M15").Value = -2025.91382

$ws.Range("H100").Value = 60533.41
$ws.Range("I100").Value = 73040.64
$ws.Range("K100").Value = 73040.64
$ws.Range("M100").Value = -72499.64

$ws.Range("H138").Value = 3651.1924
$ws.Range("I138").Value = 2356.111
$ws.Range("J138").Value = 6565.125
$ws.Range("K138").Value = 7068.333
$ws.Range("L138").Value = 19695.375
$ws.Range("M138").Value = -1928.333
$ws.Range("N138").Value = -29975.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2705.9583
$ws.Range("I32").Value = 1940.6097
$ws.Range("K32").Value = 1940.6097
$ws.Range("M32").Value = -1653.6097

$ws.Range("H44").Value = 15964.167
$ws.Range("I44").Value = 15948
$ws.Range("J44").Value = 15967.4
$ws.Range("K44").Value = 15948
$ws.Range("L44").Value = 15967.4
$ws.Range("M44").Value = -15460
$ws.Range("N44").Value = -16943.4

$ws.Range("H45").Value = 3281.2222
$ws.Range("I45").Value = 3568.875
$ws.Range("K45").Value = 3568.875
$ws.Range("M45").Value = -3191.875

$ws.Range("H61").Value = 5259.875
$ws.Range("I61").Value = 2430
$ws.Range("J61").Value = 13749.5
$ws.Range("K61").Value = 2430
$ws.Range("L61").Value = 13749.5
$ws.Range("M61").Value = -2218
$ws.Range("N61").Value = -14173.5

$ws.Range("H74").Value = 4028.2727
$ws.Range("I74").Value = 2130.2727
$ws.Range("J74").Value = 5926.273
$ws.Range("K74").Value = 2130.2727
$ws.Range("L74").Value = 5926.273
$ws.Range("M74").Value = -1256.2727
$ws.Range("N74").Value = -7674.273

$ws.Range("H77").Value = 4028.2727
$ws.Range("I77").Value = 2130.2727
$ws.Range("J77").Value = 5926.273
$ws.Range("K77").Value = 10651.3635
$ws.Range("L77").Value = 29631.365
$ws.Range("M77").Value = -6283.363499999999
$ws.Range("N77").Value = -38367.36500000001

$ws.Range("H132").Value = 6200.9414
$ws.Range("I132").Value = 6146.8477
$ws.Range("K132").Value = 18440.5431
$ws.Range("M132").Value = -15910.5431

$ws.Range("H136").Value = 5259.875
$ws.Range("I136").Value = 2430
$ws.Range("J136").Value = 13749.5
$ws.Range("K136").Value = 7290
$ws.Range("L136").Value = 41248.5
$ws.Range("M136").Value = -4740
$ws.Range("N136").Value = -46348.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 205231.5
$ws.Range("J70").Value = 205231.5
$ws.Range("L70").Value = 205231.5
$ws.Range("N70").Value = -205817.5

$ws.Range("H73").Value = 205231.5
$ws.Range("J73").Value = 205231.5
$ws.Range("L73").Value = 205231.5
$ws.Range("N73").Value = -207259.5

$ws.Range("H134").Value = 5976.548
$ws.Range("I134").Value = 3868.7097
$ws.Range("K134").Value = 11606.1291
$ws.Range("M134").Value = -9071.1291

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4610.737
$ws.Range("I58").Value = 2139.5557
$ws.Range("J58").Value = 6834.8
$ws.Range("K58").Value = 2139.5557
$ws.Range("L58").Value = 6834.8
$ws.Range("M58").Value = -1936.5557
$ws.Range("N58").Value = -7240.8

$ws.Range("H105").Value = 1641.0769
$ws.Range("I105").Value = 1639.7
$ws.Range("K105").Value = 1639.7
$ws.Range("M105").Value = 107.3

$ws.Range("H136").Value = 4610.737
$ws.Range("I136").Value = 2139.5557
$ws.Range("J136").Value = 6834.8
$ws.Range("K136").Value = 6418.6671
$ws.Range("L136").Value = 20504.4
$ws.Range("M136").Value = -3868.6671
$ws.Range("N136").Value = -25604.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4353
$ws.Range("I3").Value = 4353
$ws.Range("K3").Value = 13059
$ws.Range("M3").Value = -12947

$ws.Range("H5").Value = 1674.3478
$ws.Range("I5").Value = 651.8333
$ws.Range("K5").Value = 1955.4999
$ws.Range("M5").Value = -1843.4999

$ws.Range("H14").Value = 3190
$ws.Range("I14").Value = 3190
$ws.Range("K14").Value = 9570
$ws.Range("M14").Value = -9397

$ws.Range("H68").Value = 1161.6666
$ws.Range("J68").Value = 993
$ws.Range("L68").Value = 2979
$ws.Range("N68").Value = -4601

$ws.Range("H71").Value = 1161.6666
$ws.Range("J71").Value = 993
$ws.Range("L71").Value = 8937
$ws.Range("N71").Value = -17049

$ws.Range("H107").Value = 830.5833
$ws.Range("J107").Value = 906.375
$ws.Range("L107").Value = 2719.125
$ws.Range("N107").Value = -6559.125

$ws.Range("H113").Value = 1114.1111
$ws.Range("I113").Value = 1561
$ws.Range("J113").Value = 890.6667
$ws.Range("K113").Value = 4683
$ws.Range("L113").Value = 2672.0001
$ws.Range("M113").Value = -2513
$ws.Range("N113").Value = -7012.0001

$ws.Range("H135").Value = 1674.3478
$ws.Range("I135").Value = 651.8333
$ws.Range("K135").Value = 5866.4997
$ws.Range("M135").Value = -3331.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 37571.223
$ws.Range("J57").Value = 47147.285
$ws.Range("L57").Value = 47147.285
$ws.Range("N57").Value = -48787.285

$ws.Range("H80").Value = 3055.8
$ws.Range("I80").Value = 2994.5715
$ws.Range("K80").Value = 2994.5715
$ws.Range("M80").Value = -1996.5715

$ws.Range("H83").Value = 3055.8
$ws.Range("I83").Value = 2994.5715
$ws.Range("K83").Value = 14972.8575
$ws.Range("M83").Value = -9980.8575

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H126").Value = 7497
$ws.Range("J126").Value = 7497
$ws.Range("L126").Value = 22491
$ws.Range("N126").Value = -27431

$ws.Range("H136").Value = 37967.43
$ws.Range("J136").Value = 37967.43
$ws.Range("L136").Value = 113902.29
$ws.Range("N136").Value = -119002.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2907.818
$ws.Range("I22").Value = 2272.926
$ws.Range("J22").Value = 5764.8335
$ws.Range("K22").Value = 2272.926
$ws.Range("L22").Value = 5764.8335
$ws.Range("M22").Value = -1977.926
$ws.Range("N22").Value = -6354.8335

$ws.Range("H27").Value = 2907.818
$ws.Range("I27").Value = 2272.926
$ws.Range("J27").Value = 5764.8335
$ws.Range("K27").Value = 2272.926
$ws.Range("L27").Value = 5764.8335
$ws.Range("M27").Value = -2165.926
$ws.Range("N27").Value = -5978.8335

$ws.Range("H40").Value = 2999.6667
$ws.Range("I40").Value = 2999.5
$ws.Range("K40").Value = 2999.5
$ws.Range("M40").Value = -2863.5

$ws.Range("H46").Value = 1444.5
$ws.Range("I46").Value = 998.8
$ws.Range("J46").Value = 1815.9166
$ws.Range("K46").Value = 998.8
$ws.Range("L46").Value = 1815.9166
$ws.Range("M46").Value = -810.8
$ws.Range("N46").Value = -2191.9166

$ws.Range("H55").Value = 1532.1666
$ws.Range("I55").Value = 1981
$ws.Range("J55").Value = 1083.3334
$ws.Range("K55").Value = 1981
$ws.Range("L55").Value = 1083.3334
$ws.Range("M55").Value = -1808
$ws.Range("N55").Value = -1429.3334

$ws.Range("H82").Value = 2403.2144
$ws.Range("I82").Value = 1799.8889
$ws.Range("K82").Value = 1799.8889
$ws.Range("M82").Value = -1438.8889

$ws.Range("H85").Value = 2403.2144
$ws.Range("I85").Value = 1799.8889
$ws.Range("K85").Value = 1799.8889
$ws.Range("M85").Value = -551.8888999999999

$ws.Range("H93").Value = 1705.7727
$ws.Range("I93").Value = 1740.4474
$ws.Range("J93").Value = 1486.1666
$ws.Range("K93").Value = 1740.4474
$ws.Range("L93").Value = 1486.1666
$ws.Range("M93").Value = -492.4474
$ws.Range("N93").Value = -3982.1666

$ws.Range("H100").Value = 1213911.6
$ws.Range("I100").Value = 128125.125
$ws.Range("K100").Value = 128125.125
$ws.Range("M100").Value = -127584.125

$ws.Range("H122").Value = 3137.8572
$ws.Range("I122").Value = 2804.111
$ws.Range("J122").Value = 3738.6
$ws.Range("K122").Value = 8412.332999999999
$ws.Range("L122").Value = 11215.8
$ws.Range("M122").Value = -5962.332999999999
$ws.Range("N122").Value = -16115.8

$ws.Range("H132").Value = 6177867.5
$ws.Range("I132").Value = 7941116
$ws.Range("K132").Value = 23823348
$ws.Range("M132").Value = -23820818

$ws.Range("H136").Value = 13894722
$ws.Range("J136").Value = 15000
$ws.Range("L136").Value = 45000
$ws.Range("N136").Value = -50100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5135.727
$ws.Range("J122").Value = 4390.25
$ws.Range("L122").Value = 13170.75
$ws.Range("N122").Value = -18070.75

$ws.Range("H132").Value = 15556.015
$ws.Range("I132").Value = 9441.611000000001
$ws.Range("K132").Value = 28324.833
$ws.Range("M132").Value = -25794.833

